$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 5000
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 5000
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 5000
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -5586
$ws.Range("H11").Value = 1515
$ws.Range("I11").Value = 1515
$ws.Range("K11").Value = 1515
$ws.Range("M11").Value = -1375
$ws.Range("H15").Value = 329.9375
$ws.Range("I15").Value = 329.9375
$ws.Range("K15").Value = 989.8125
$ws.Range("M15").Value = -820.8125
$ws.Range("H38").Value = 30
$ws.Range("I38").Value = 30
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 90
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 282
$ws.Range("N38").ClearContents()
$ws.Range("H39").Value = 30.88889
$ws.Range("I39").Value = 41.333332
$ws.Range("J39").Value = 10
$ws.Range("K39").Value = 123.999996
$ws.Range("L39").Value = 30
$ws.Range("M39").Value = 172.000004
$ws.Range("N39").Value = -622

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 2173
$ws.Range("I41").Value = 2173
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 2173
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -1759
$ws.Range("N41").ClearContents()
$ws.Range("H63").Value = 7494
$ws.Range("J63").Value = 7494
$ws.Range("L63").Value = 7494
$ws.Range("N63").Value = -8866
$ws.Range("H66").Value = 7494
$ws.Range("J66").Value = 7494
$ws.Range("L66").Value = 37470
$ws.Range("N66").Value = -44334
$ws.Range("H92").Value = 97666.336
$ws.Range("J92").Value = 97666.336
$ws.Range("L92").Value = 97666.336
$ws.Range("N92").Value = -102658.336
$ws.Range("H122").Value = 2603
$ws.Range("I122").Value = 2314.3
$ws.Range("J122").Value = 5490
$ws.Range("K122").Value = 6942.900000000001
$ws.Range("L122").Value = 16470
$ws.Range("M122").Value = -4492.900000000001
$ws.Range("N122").Value = -21370

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H115").Value = 65000
$ws.Range("J115").Value = 65000
$ws.Range("L115").Value = 65000
$ws.Range("N115").Value = -68134
$ws.Range("H134").Value = 7956.2
$ws.Range("I134").Value = 1794.7142
$ws.Range("K134").Value = 5384.142599999999
$ws.Range("M134").Value = -2849.142599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H57").Value = 3833.3333
$ws.Range("I57").Value = 3833.3333
$ws.Range("K57").Value = 3833.3333
$ws.Range("M57").Value = -3273.3333
$ws.Range("H99").Value = 4605.6
$ws.Range("I99").Value = 4605.6
$ws.Range("K99").Value = 4605.6
$ws.Range("M99").Value = -3107.6
$ws.Range("H126").Value = 4605.6
$ws.Range("I126").Value = 4605.6
$ws.Range("K126").Value = 13816.8
$ws.Range("M126").Value = -11346.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 230
$ws.Range("J15").Value = 262.7143
$ws.Range("L15").Value = 788.1428999999999
$ws.Range("N15").Value = -1068.1429
$ws.Range("H26").Value = 750
$ws.Range("I26").Value = 750
$ws.Range("K26").Value = 2250
$ws.Range("M26").Value = -1962
$ws.Range("H36").Value = 3745
$ws.Range("J36").Value = 3745
$ws.Range("L36").Value = 11235
$ws.Range("N36").Value = -11573
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H95").Value = 38206
$ws.Range("J95").Value = 38206
$ws.Range("L95").Value = 38206
$ws.Range("N95").Value = -43698
$ws.Range("H122").Value = 7296.5
$ws.Range("I122").Value = 6002.3335
$ws.Range("K122").Value = 18007.0005
$ws.Range("M122").Value = -15557.0005
$ws.Range("H126").Value = 1421.5
$ws.Range("J126").Value = 1399
$ws.Range("L126").Value = 4197
$ws.Range("N126").Value = -9137

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H55").Value = 869.6
$ws.Range("I55").Value = 892.75
$ws.Range("J55").Value = 777
$ws.Range("K55").Value = 892.75
$ws.Range("L55").Value = 777
$ws.Range("M55").Value = -719.75
$ws.Range("N55").Value = -1123
$ws.Range("H74").Value = 24996.666
$ws.Range("I74").Value = 24990
$ws.Range("K74").Value = 24990
$ws.Range("M74").Value = -23992
$ws.Range("H77").Value = 24996.666
$ws.Range("I77").Value = 24990
$ws.Range("K77").Value = 74970
$ws.Range("M77").Value = -69978
$ws.Range("H109").Value = 30000
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("H132").Value = 10549.25
$ws.Range("I132").Value = 4879
$ws.Range("K132").Value = 14637
$ws.Range("M132").Value = -12107
$ws.Range("H136").Value = 16166.5
$ws.Range("I136").Value = 12666.333
$ws.Range("J136").Value = 19666.666
$ws.Range("K136").Value = 37998.999
$ws.Range("L136").Value = 58999.99800000001
$ws.Range("M136").Value = -35448.999
$ws.Range("N136").Value = -64099.99800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 24997.5
$ws.Range("J70").Value = 25000
$ws.Range("L70").Value = 25000
$ws.Range("N70").Value = -25630
$ws.Range("H73").Value = 24997.5
$ws.Range("J73").Value = 25000
$ws.Range("L73").Value = 25000
$ws.Range("N73").Value = -27184
$ws.Range("H98").Value = 19999
$ws.Range("J98").Value = 19999
$ws.Range("L98").Value = 19999
$ws.Range("N98").Value = -25989
$ws.Range("H132").Value = 7621.6113
$ws.Range("I132").Value = 6346
$ws.Range("J132").Value = 13999.667
$ws.Range("K132").Value = 19038
$ws.Range("L132").Value = 41999.001
$ws.Range("M132").Value = -16508
$ws.Range("N132").Value = -47059.001
$ws.Range("H136").Value = 5762.6
$ws.Range("I136").Value = 2356.2
$ws.Range("K136").Value = 7068.599999999999
$ws.Range("M136").Value = -4518.599999999999
